$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Workbook-wide status text: every cell currently showing
#    "Ready for handoff" becomes "Handed back: in sync with en-US"
#    (Overview!B2:C3 and the "Status" column (C) on each language sheet).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet ("Generate Report for Handback"): fill in the
#    "Latest Target File" (F) / "Latest Handback File" (G) columns and bump
#    the "Latest Handback DateTime" (H) for both data rows.
# ---------------------------------------------------------------------------
$zhTargetFile = "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
$zhHandbackFile = "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
$zhTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/dc57b49458b1653baf2bf997359eaecfb28cb93c/e2e/144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46888786891e2c36332fefad66d4496d0c635a9e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
$zhHandbackDate = "2016-03-21 03:42:47"

$wsZhCn.Range("F2").Value = $zhTargetFile
$wsZhCn.Range("G2").Value = $zhHandbackFile
$wsZhCn.Range("F3").Value = $zhTargetFile
$wsZhCn.Range("G3").Value = $zhHandbackFile

$wsZhCn.Range("H2").Value = $zhHandbackDate
$wsZhCn.Range("H3").Value = $zhHandbackDate

# Rebuild the hyperlinks in the final left-to-right / top-to-bottom order so
# relationship ids line up the way Excel would renumber them: A2, D2, F2, G2,
# A3, D3, F3, G3.
$zhA2Url = "https://github.com/OpenLocalizationTest/oltest/blob/dc57b49458b1653baf2bf997359eaecfb28cb93c/e2e/144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
$zhD2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46888786891e2c36332fefad66d4496d0c635a9e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
$zhA3Url = "https://github.com/OpenLocalizationTest/oltest/blob/dc57b49458b1653baf2bf997359eaecfb28cb93c/e2e/ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
$zhD3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46888786891e2c36332fefad66d4496d0c635a9e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhA2Url, "", "", "144dc3ed-8811-48a6-98eb-0053a6c9080a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhD2Url, "", "", $zhHandbackFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhTargetUrl, "", "", $zhTargetFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhHandbackUrl, "", "", $zhHandbackFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhA3Url, "", "", "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhD3Url, "", "", $zhHandbackFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhTargetUrl, "", "", $zhTargetFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhHandbackUrl, "", "", $zhHandbackFile)

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of change, but the handback datetime + xlf file
#    name are specific to this locale.
# ---------------------------------------------------------------------------
$deTargetFile = "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
$deHandbackFile = "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
$deTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/dc57b49458b1653baf2bf997359eaecfb28cb93c/e2e/144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ed2cfcf5b37b8a73d4d857b5e72b7aff74c0928/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
$deHandbackDate = "2016-03-21 03:43:00"

$wsDeDe.Range("F2").Value = $deTargetFile
$wsDeDe.Range("G2").Value = $deHandbackFile
$wsDeDe.Range("F3").Value = $deTargetFile
$wsDeDe.Range("G3").Value = $deHandbackFile

$wsDeDe.Range("H2").Value = $deHandbackDate
$wsDeDe.Range("H3").Value = $deHandbackDate

$deA2Url = "https://github.com/OpenLocalizationTest/oltest/blob/dc57b49458b1653baf2bf997359eaecfb28cb93c/e2e/144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
$deD2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ed2cfcf5b37b8a73d4d857b5e72b7aff74c0928/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
$deA3Url = "https://github.com/OpenLocalizationTest/oltest/blob/dc57b49458b1653baf2bf997359eaecfb28cb93c/e2e/ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
$deD3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ed2cfcf5b37b8a73d4d857b5e72b7aff74c0928/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deA2Url, "", "", "144dc3ed-8811-48a6-98eb-0053a6c9080a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deD2Url, "", "", $deHandbackFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deTargetUrl, "", "", $deTargetFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deHandbackUrl, "", "", $deHandbackFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deA3Url, "", "", "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deD3Url, "", "", $deHandbackFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deTargetUrl, "", "", $deTargetFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deHandbackUrl, "", "", $deHandbackFile)

Write-Host "Report generated for handback."
